$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-12 (column B - count) with revised values
$ws.Range("B2").Value = 78477
$ws.Range("B3").Value = 73442
$ws.Range("B4").Value = 67464
$ws.Range("B5").Value = 61710
$ws.Range("B6").Value = 73578
$ws.Range("B7").Value = 126747
$ws.Range("B8").Value = 92712
$ws.Range("B9").Value = 83152
$ws.Range("B10").Value = 93596
$ws.Range("B11").Value = 81302
$ws.Range("B12").Value = 75451

# Add new rows 13-15 for the additional quarters
$ws.Range("A13").Value = 44561
$ws.Range("B13").Value = 64443
$ws.Range("C13").Value = "Q"

$ws.Range("A14").Value = 44651
$ws.Range("B14").Value = 67004
$ws.Range("C14").Value = "Q"

$ws.Range("A15").Value = 44742
$ws.Range("B15").Value = 45389
$ws.Range("C15").Value = "Q"

# Apply the same date number format used by the existing ts column cells
$ws.Range("A13:A15").NumberFormat = $ws.Range("A12").NumberFormat
